$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L: header "break_on_off"
$ws.Range("L1").Value = "break_on_off"

# Fill data rows 2-73 with 0, except the break rows which get 1
$breakRows = @(19,37,54)
for ($r = 2; $r -le 73; $r++) {
    if ($breakRows -contains $r) {
        $ws.Cells.Item($r, 12).Value = 1
    } else {
        $ws.Cells.Item($r, 12).Value = 0
    }
}

# Update selection to match the saved workbook view
$null = $ws.Range("A2").Select()
